$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5326
$ws.Range("J3").Value = 5652
$ws.Range("C4").Value = 1834
$ws.Range("D4").Value = 1961
$ws.Range("F4").Value = 1899
$ws.Range("J4").Value = 1250
$ws.Range("J5").Value = 441
$ws.Range("J6").Value = 7082
$ws.Range("C7").Value = 28378
$ws.Range("D7").Value = 28151
$ws.Range("F7").Value = 24090
$ws.Range("J7").Value = 19751

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 76
$ws.Range("J8").Value = 1247
$ws.Range("J9").Value = 98
$ws.Range("J15").Value = 217
$ws.Range("J19").Value = 568
$ws.Range("J20").Value = 411
$ws.Range("J24").Value = 62
$ws.Range("J25").Value = 99
$ws.Range("J29").Value = 1103
$ws.Range("J33").Value = 911
$ws.Range("J34").Value = 93
$ws.Range("J36").Value = 274
$ws.Range("J37").Value = 615
$ws.Range("J42").Value = 804
$ws.Range("J43").Value = 165
$ws.Range("J47").Value = 148
$ws.Range("J48").Value = 228
$ws.Range("J49").Value = 133
$ws.Range("J50").Value = 123
$ws.Range("J51").Value = 248
$ws.Range("J52").Value = 499
$ws.Range("J53").Value = 267
$ws.Range("J55").Value = 257
$ws.Range("J57").Value = 84
$ws.Range("J61").Value = 22
$ws.Range("C63").Value = 265
$ws.Range("D63").Value = 345
$ws.Range("F63").Value = 187
$ws.Range("J63").Value = 65
$ws.Range("J64").Value = 134
$ws.Range("J65").Value = 507
$ws.Range("J71").Value = 69
$ws.Range("J72").Value = 76
$ws.Range("J75").Value = 59
$ws.Range("J76").Value = 284
$ws.Range("J77").Value = 156
$ws.Range("J83").Value = 408
$ws.Range("J84").Value = 170
$ws.Range("J85").Value = 838
$ws.Range("J87").Value = 68
$ws.Range("J88").Value = 217
$ws.Range("J89").Value = 256
$ws.Range("J90").Value = 215
$ws.Range("J91").Value = 221
$ws.Range("J94").Value = 196
$ws.Range("J96").Value = 236
$ws.Range("J98").Value = 137
$ws.Range("J99").Value = 309
$ws.Range("C101").Value = 28378
$ws.Range("D101").Value = 28151
$ws.Range("F101").Value = 24090
$ws.Range("J101").Value = 19751

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 82
$ws.Range("J3").Value = 72
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 256

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 217
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 241
$ws.Range("J7").Value = 838

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 148
$ws.Range("J7").Value = 499

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 52
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 169
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 345
$ws.Range("J6").Value = 415
$ws.Range("J7").Value = 1247

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 152
$ws.Range("J4").Value = 11
$ws.Range("J5").Value = 13
$ws.Range("J7").Value = 408

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 223
$ws.Range("J3").Value = 300
$ws.Range("J6").Value = 310
$ws.Range("J7").Value = 911

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 212
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 615

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 181
$ws.Range("J7").Value = 507

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J5").Value = 9
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J4").Value = 62
$ws.Range("J6").Value = 291
$ws.Range("J7").Value = 1103

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 142
$ws.Range("J3").Value = 168
$ws.Range("J7").Value = 568

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 58
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 182
$ws.Range("J6").Value = 408
$ws.Range("J7").Value = 804

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 67
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 142
$ws.Range("J7").Value = 411

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 89
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 274

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J2").Value = 27
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 61
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 58
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 156

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 22
